$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sending cluster label in row 2 (A2): "Inflammatory-Mac" -> "ECs"
$ws.Range("A2").Value = "ECs"

# Row 2 updated metric values
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5679715
$ws.Range("H2").Value = 1.135943
$ws.Range("I2").Value = 0.1948020679094191
$ws.Range("J2").Value = 0.1388865047139418
$ws.Range("M2").Value = 0.115994
$ws.Range("N2").Value = 0.231988
$ws.Range("Q2").Value = 0.06588128617099999
$ws.Range("R2").Value = 0.263525144684
$ws.Range("S2").Value = 0.1948020679094191
$ws.Range("T2").Value = 0.1388865047139418

# Row 3 updated metric values
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.235202
$ws.Range("H3").Value = 6.705606
$ws.Range("I3").Value = 0.7666264448044829
$ws.Range("J3").Value = 0.8198634784745682
$ws.Range("M3").Value = 0.115994
$ws.Range("N3").Value = 0.231988
$ws.Range("Q3").Value = 0.259270020788
$ws.Range("R3").Value = 1.555620124728
$ws.Range("S3").Value = 0.7666264448044829
$ws.Range("T3").Value = 0.8198634784745682

# Row 4 updated metric values
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1124603333333333
$ws.Range("H4").Value = 0.337381
$ws.Range("I4").Value = 0.03857148728609783
$ws.Range("J4").Value = 0.04125001681149001
$ws.Range("M4").Value = 0.115994
$ws.Range("N4").Value = 0.231988
$ws.Range("Q4").Value = 0.01304472390466667
$ws.Range("R4").Value = 0.078268343428
$ws.Range("S4").Value = 0.03857148728609783
$ws.Range("T4").Value = 0.04125001681149001
